# Auto-generated edit script: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H, I, J, K, L, M, N) across several Leve rows in the Asura_Profits
# workbook, reflecting a refreshed market-price snapshot from the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 552.1875
$ws.Range("J17").Value = 406.67743
$ws.Range("L17").Value = 1220.03229
$ws.Range("N17").Value = -1556.03229
$ws.Range("H46").Value = 6499.8887
$ws.Range("J46").Value = 6499.8887
$ws.Range("L46").Value = 19499.6661
$ws.Range("N46").Value = -19737.6661
$ws.Range("H60").Value = 6499.8887
$ws.Range("J60").Value = 6499.8887
$ws.Range("L60").Value = 19499.6661
$ws.Range("N60").Value = -20467.6661
$ws.Range("H113").Value = 3400.8333
$ws.Range("I113").Value = 2752.5
$ws.Range("J113").Value = 3725
$ws.Range("K113").Value = 2752.5
$ws.Range("L113").Value = 3725
$ws.Range("M113").Value = 501.5
$ws.Range("N113").Value = -10233
$ws.Range("H129").Value = 1081.2764
$ws.Range("J129").Value = 1134.3715
$ws.Range("L129").Value = 3403.1145
$ws.Range("N129").Value = -13403.1145
$ws.Range("H132").Value = 2403.375
$ws.Range("I132").Value = 2290.7273
$ws.Range("J132").Value = 2651.2
$ws.Range("K132").Value = 6872.1819
$ws.Range("L132").Value = 7953.599999999999
$ws.Range("M132").Value = -4342.1819
$ws.Range("N132").Value = -13013.6
$ws.Range("H133").Value = 71525
$ws.Range("J133").Value = 71525
$ws.Range("L133").Value = 71525
$ws.Range("N133").Value = -81645
$ws.Range("H137").Value = 1626.6769
$ws.Range("I137").Value = 1002.97144
$ws.Range("J137").Value = 2354.3333
$ws.Range("K137").Value = 3008.91432
$ws.Range("L137").Value = 7062.999899999999
$ws.Range("M137").Value = -458.9143199999999
$ws.Range("N137").Value = -12162.9999
$ws.Range("H138").Value = 3795.9023
$ws.Range("I138").Value = 2367
$ws.Range("J138").Value = 4287.82
$ws.Range("K138").Value = 7101
$ws.Range("L138").Value = 12863.46
$ws.Range("M138").Value = -1961
$ws.Range("N138").Value = -23143.46

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1775.75
$ws.Range("I74").Value = 1750.5
$ws.Range("J74").Value = 1841.4
$ws.Range("K74").Value = 1750.5
$ws.Range("L74").Value = 1841.4
$ws.Range("M74").Value = -876.5
$ws.Range("N74").Value = -3589.4
$ws.Range("H77").Value = 1775.75
$ws.Range("I77").Value = 1750.5
$ws.Range("J77").Value = 1841.4
$ws.Range("K77").Value = 8752.5
$ws.Range("L77").Value = 9207
$ws.Range("M77").Value = -4384.5
$ws.Range("N77").Value = -17943
$ws.Range("H88").Value = 2326.375
$ws.Range("I88").Value = 2003.6666
$ws.Range("J88").Value = 2520
$ws.Range("K88").Value = 2003.6666
$ws.Range("L88").Value = 2520
$ws.Range("M88").Value = -1597.6666
$ws.Range("N88").Value = -3332
$ws.Range("H91").Value = 2326.375
$ws.Range("I91").Value = 2003.6666
$ws.Range("J91").Value = 2520
$ws.Range("K91").Value = 2003.6666
$ws.Range("L91").Value = 2520
$ws.Range("M91").Value = -599.6666
$ws.Range("N91").Value = -5328
$ws.Range("H97").Value = 1386.25
$ws.Range("I97").Value = 1098.3334
$ws.Range("K97").Value = 1098.3334
$ws.Range("M97").Value = -602.3334
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("H123").Value = 25427.777
$ws.Range("J123").Value = 25427.777
$ws.Range("L123").Value = 25427.777
$ws.Range("N123").Value = -35227.777
$ws.Range("H131").Value = 43290.418
$ws.Range("J131").Value = 43290.418
$ws.Range("L131").Value = 43290.418
$ws.Range("N131").Value = -53370.418

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 363.76
$ws.Range("I22").Value = 218.7619
$ws.Range("J22").Value = 1125
$ws.Range("K22").Value = 218.7619
$ws.Range("L22").Value = 1125
$ws.Range("M22").Value = 131.2381
$ws.Range("N22").Value = -1825
$ws.Range("H31").Value = 1546.7174
$ws.Range("I31").Value = 1726.5454
$ws.Range("J31").Value = 1381.875
$ws.Range("K31").Value = 1726.5454
$ws.Range("L31").Value = 1381.875
$ws.Range("M31").Value = -1431.5454
$ws.Range("N31").Value = -1971.875
$ws.Range("H34").Value = 1546.7174
$ws.Range("I34").Value = 1726.5454
$ws.Range("J34").Value = 1381.875
$ws.Range("K34").Value = 1726.5454
$ws.Range("L34").Value = 1381.875
$ws.Range("M34").Value = -1524.5454
$ws.Range("N34").Value = -1785.875
$ws.Range("H107").Value = 1230.5
$ws.Range("I107").Value = 1230.5
$ws.Range("K107").Value = 1230.5
$ws.Range("M107").Value = 689.5
$ws.Range("H133").Value = 49672
$ws.Range("J133").Value = 49672
$ws.Range("L133").Value = 49672
$ws.Range("N133").Value = -54732

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1355
$ws.Range("J29").Value = 1226
$ws.Range("L29").Value = 3678
$ws.Range("N29").Value = -4232
$ws.Range("H34").Value = 1229.2858
$ws.Range("J34").Value = 3140
$ws.Range("L34").Value = 9420
$ws.Range("N34").Value = -9588
$ws.Range("H39").Value = 4189
$ws.Range("J39").Value = 4832.4546
$ws.Range("L39").Value = 14497.3638
$ws.Range("N39").Value = -15085.3638
$ws.Range("H68").Value = 1273.6
$ws.Range("I68").Value = 940.0714
$ws.Range("K68").Value = 2820.2142
$ws.Range("M68").Value = -2009.2142
$ws.Range("H71").Value = 1273.6
$ws.Range("I71").Value = 940.0714
$ws.Range("K71").Value = 8460.642600000001
$ws.Range("M71").Value = -4404.642600000001
$ws.Range("H97").Value = 1822.8
$ws.Range("J97").Value = 1822.8
$ws.Range("L97").Value = 5468.4
$ws.Range("N97").Value = -6460.4
$ws.Range("H102").Value = 7999.4165
$ws.Range("J102").Value = 7999.4165
$ws.Range("L102").Value = 23998.2495
$ws.Range("N102").Value = -28866.2495
$ws.Range("H107").Value = 1339.0461
$ws.Range("I107").Value = 1211.6923
$ws.Range("J107").Value = 1530.0769
$ws.Range("K107").Value = 3635.0769
$ws.Range("L107").Value = 4590.2307
$ws.Range("M107").Value = -1715.0769
$ws.Range("N107").Value = -8430.2307
$ws.Range("H132").Value = 1347.0344
$ws.Range("J132").Value = 1402.826
$ws.Range("L132").Value = 12625.434
$ws.Range("N132").Value = -17685.434
$ws.Range("H139").Value = 2349.3333
$ws.Range("I139").Value = 2219.2
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 6657.599999999999
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -1517.599999999999
$ws.Range("N139").Value = -19280
$ws.Range("H140").Value = 1652.1562
$ws.Range("I140").Value = 822.8095
$ws.Range("J140").Value = 3235.4546
$ws.Range("K140").Value = 2468.4285
$ws.Range("L140").Value = 9706.363799999999
$ws.Range("M140").Value = 2711.5715
$ws.Range("N140").Value = -20066.3638

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8604.647000000001
$ws.Range("J123").Value = 8604.647000000001
$ws.Range("L123").Value = 8604.647000000001
$ws.Range("N123").Value = -13504.647

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2333.3333
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 2333.3333
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H75").Value = 26000
$ws.Range("J75").Value = 26000
$ws.Range("L75").Value = 26000
$ws.Range("N75").Value = -27872
$ws.Range("H78").Value = 26000
$ws.Range("J78").Value = 26000
$ws.Range("L78").Value = 78000
$ws.Range("N78").Value = -87360
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2300
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 1600
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 1600
$ws.Range("M96").Value = -1627
$ws.Range("N96").Value = -4346
$ws.Range("H123").Value = 22875.414
$ws.Range("J123").Value = 22875.414
$ws.Range("L123").Value = 22875.414
$ws.Range("N123").Value = -32675.414
$ws.Range("H132").Value = 2146.3845
$ws.Range("I132").Value = 2031.3462
$ws.Range("K132").Value = 6094.0386
$ws.Range("M132").Value = -3564.0386

